$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price cells stay as text (matches source formatting)
$textCells = @("D5", "D6", "D8", "D10", "D12", "D13", "D18", "D19", "D21", "D23", "D24", "D25", "D27", "D28", "D29", "D31", "D33", "D36", "D37", "D38", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range('D2').Value = '62.549.35'
$ws.Range('E2').Value = '  -1.06%  '
$ws.Range('D3').Value = '2.562.72'
$ws.Range('E3').Value = '  -0.25%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '578.29'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('D6').Value = '143.22'
$ws.Range('E6').Value = '  -4.04%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('D8').Value = '0.586'
$ws.Range('E8').Value = '  -0.17%  '
$ws.Range('E9').Value = '  -3.29%  '
$ws.Range('D10').Value = '5.53'
$ws.Range('E10').Value = '  -1.23%  '
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = '0.350'
$ws.Range('E12').Value = '  -2.25%  '
$ws.Range('D13').Value = '26.74'
$ws.Range('E13').Value = '  -4.44%  '
$ws.Range('D14').Value = '3.022.35'
$ws.Range('E14').Value = '  -0.07%  '
$ws.Range('D15').Value = '62.508.36'
$ws.Range('E15').Value = '  -0.89%  '
$ws.Range('E16').Value = '  -3.33%  '
$ws.Range('D17').Value = '2.563.50'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').Value = '11.09'
$ws.Range('E18').Value = '  -3.25%  '
$ws.Range('D19').Value = '337.55'
$ws.Range('E19').Value = '  -1.34%  '
$ws.Range('E20').Value = '  -2.01%  '
$ws.Range('D21').Value = '6.63'
$ws.Range('E21').Value = '  -3.32%  '
$ws.Range('E22').Value = '  -0.02%  '
$ws.Range('D23').Value = '67.09'
$ws.Range('E23').Value = '  +1.41%  '
$ws.Range('B24').Value = 'Fetch.AI'
$ws.Range('C24').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D24').Value = '1.57'
$ws.Range('E24').Value = '  -5.05%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').Value = '0.163'
$ws.Range('E25').Value = '  -4.26%  '
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('D28').Value = '7.90'
$ws.Range('E28').Value = '  -5.21%  '
$ws.Range('D29').Value = '8.16'
$ws.Range('E29').Value = '  -4.99%  '
$ws.Range('E30').Value = '  -2.68%  '
$ws.Range('D31').Value = '454.08'
$ws.Range('E31').Value = '  +3.16%  '
$ws.Range('D32').Value = '0.0₃0795'
$ws.Range('E32').Value = '  -4.65%  '
$ws.Range('D33').Value = '176.60'
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('D36').Value = '0.394'
$ws.Range('E36').Value = '  -3.19%  '
$ws.Range('D37').Value = '18.78'
$ws.Range('E37').Value = '  -3.07%  '
$ws.Range('D38').Value = '4.42'
$ws.Range('E38').Value = '  -2.98%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('E40').Value = '  -4.77%  '
$ws.Range('B41').Value = 'OKB'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D41').Value = '40.16'
$ws.Range('E41').Value = '  +0.98%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').Value = '157.74'
$ws.Range('E42').Value = '  +3.51%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = '3.67'
$ws.Range('E43').Value = '  -4.06%  '
$ws.Range('B44').Value = 'Mantle'
$ws.Range('C44').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value = '0.627'
$ws.Range('E44').Value = '  +3.01%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '20.79'
$ws.Range('E45').Value = '  -3.69%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '0.0530'
$ws.Range('E46').Value = '  -5.13%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '0.0956'
$ws.Range('E47').Value = '  -2.39%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0233'
$ws.Range('E48').Value = '  -4.02%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '17.92'
$ws.Range('E49').Value = '  -3.29%  '
$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').Value = '11.40'
$ws.Range('E50').Value = '  +0.18%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = '1.66'
$ws.Range('E51').Value = '  -5.87%  '
